$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-06 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-07 Sunday", 2) | Out-Null
$d.Content.Find.Execute("36-27=9", $true, $false, $false, $false, $false, $true, 1, $false, "35+58=93", 2) | Out-Null
$d.Content.Find.Execute("80-1=79", $true, $false, $false, $false, $false, $true, 1, $false, "61-12=49", 2) | Out-Null
$d.Content.Find.Execute("29+29=58", $true, $false, $false, $false, $false, $true, 1, $false, "9+7=16", 2) | Out-Null
$d.Content.Find.Execute("28+48=76", $true, $false, $false, $false, $false, $true, 1, $false, "44-38=6", 2) | Out-Null
$d.Content.Find.Execute("37+18=55", $true, $false, $false, $false, $false, $true, 1, $false, "83+8=91", 2) | Out-Null
$d.Content.Find.Execute("44-15=29", $true, $false, $false, $false, $false, $true, 1, $false, "65-49=16", 2) | Out-Null
$d.Content.Find.Execute("90-63=27", $true, $false, $false, $false, $false, $true, 1, $false, "35+16=51", 2) | Out-Null
$d.Content.Find.Execute("52-38=14", $true, $false, $false, $false, $false, $true, 1, $false, "46+8=54", 2) | Out-Null
$d.Content.Find.Execute("62-59=3", $true, $false, $false, $false, $false, $true, 1, $false, "63-45=18", 2) | Out-Null
$d.Content.Find.Execute("74-9=65", $true, $false, $false, $false, $false, $true, 1, $false, "7+84=91", 2) | Out-Null
$d.Content.Find.Execute("6+7=13", $true, $false, $false, $false, $false, $true, 1, $false, "59+24=83", 2) | Out-Null
$d.Content.Find.Execute("39+45=84", $true, $false, $false, $false, $false, $true, 1, $false, "49+18=67", 2) | Out-Null
$d.Content.Find.Execute("53-9=44", $true, $false, $false, $false, $false, $true, 1, $false, "90-21=69", 2) | Out-Null
$d.Content.Find.Execute("27+64=91", $true, $false, $false, $false, $false, $true, 1, $false, "79+16=95", 2) | Out-Null
$d.Content.Find.Execute("65-48=17", $true, $false, $false, $false, $false, $true, 1, $false, "63-34=29", 2) | Out-Null
$d.Content.Find.Execute("92-13=79", $true, $false, $false, $false, $false, $true, 1, $false, "56+28=84", 2) | Out-Null
$d.Content.Find.Execute("60-45=15", $true, $false, $false, $false, $false, $true, 1, $false, "68+3=71", 2) | Out-Null
$d.Content.Find.Execute("64-46=18", $true, $false, $false, $false, $false, $true, 1, $false, "44-7=37", 2) | Out-Null
$d.Content.Find.Execute("26+25=51", $true, $false, $false, $false, $false, $true, 1, $false, "39+22=61", 2) | Out-Null
$d.Content.Find.Execute("51-47=4", $true, $false, $false, $false, $false, $true, 1, $false, "87+9=96", 2) | Out-Null
$d.Content.Find.Execute("53+8=61", $true, $false, $false, $false, $false, $true, 1, $false, "92-25=67", 2) | Out-Null
$d.Content.Find.Execute("20-17=3", $true, $false, $false, $false, $false, $true, 1, $false, "14-8=6", 2) | Out-Null
$d.Content.Find.Execute("23+19=42", $true, $false, $false, $false, $false, $true, 1, $false, "7+67=74", 2) | Out-Null
$d.Content.Find.Execute("84-59=25", $true, $false, $false, $false, $false, $true, 1, $false, "13+29=42", 2) | Out-Null
$d.Content.Find.Execute("80-13=67", $true, $false, $false, $false, $false, $true, 1, $false, "9+57=66", 2) | Out-Null
$d.Content.Find.Execute("55+28=83", $true, $false, $false, $false, $false, $true, 1, $false, "32+39=71", 2) | Out-Null
$d.Content.Find.Execute("25+59=84", $true, $false, $false, $false, $false, $true, 1, $false, "72-34=38", 2) | Out-Null
$d.Content.Find.Execute("5+89=94", $true, $false, $false, $false, $false, $true, 1, $false, "16+66=82", 2) | Out-Null
$d.Content.Find.Execute("60-36=24", $true, $false, $false, $false, $false, $true, 1, $false, "39+28=67", 2) | Out-Null
$d.Content.Find.Execute("30-24=6", $true, $false, $false, $false, $false, $true, 1, $false, "80-32=48", 2) | Out-Null
$d.Content.Find.Execute("70-56=14", $true, $false, $false, $false, $false, $true, 1, $false, "18+74=92", 2) | Out-Null
$d.Content.Find.Execute("80-77=3", $true, $false, $false, $false, $false, $true, 1, $false, "66+17=83", 2) | Out-Null
$d.Content.Find.Execute("32-6=26", $true, $false, $false, $false, $false, $true, 1, $false, "93-44=49", 2) | Out-Null
$d.Content.Find.Execute("39+36=75", $true, $false, $false, $false, $false, $true, 1, $false, "64+29=93", 2) | Out-Null
$d.Content.Find.Execute("70-47=23", $true, $false, $false, $false, $false, $true, 1, $false, "37+7=44", 2) | Out-Null
$d.Content.Find.Execute("6+15=21", $true, $false, $false, $false, $false, $true, 1, $false, "81-47=34", 2) | Out-Null
$d.Content.Find.Execute("17+7=24", $true, $false, $false, $false, $false, $true, 1, $false, "91-65=26", 2) | Out-Null
$d.Content.Find.Execute("19+69=88", $true, $false, $false, $false, $false, $true, 1, $false, "82-45=37", 2) | Out-Null
$d.Content.Find.Execute("69+13=82", $true, $false, $false, $false, $false, $true, 1, $false, "39+38=77", 2) | Out-Null
$d.Content.Find.Execute("19+9=28", $true, $false, $false, $false, $false, $true, 1, $false, "70-9=61", 2) | Out-Null
$d.Content.Find.Execute("90-28=62", $true, $false, $false, $false, $false, $true, 1, $false, "77+15=92", 2) | Out-Null
$d.Content.Find.Execute("9+22=31", $true, $false, $false, $false, $false, $true, 1, $false, "94-35=59", 2) | Out-Null
$d.Content.Find.Execute("56+27=83", $true, $false, $false, $false, $false, $true, 1, $false, "18+45=63", 2) | Out-Null
$d.Content.Find.Execute("92-35=57", $true, $false, $false, $false, $false, $true, 1, $false, "48+7=55", 2) | Out-Null
$d.Content.Find.Execute("4+9=13", $true, $false, $false, $false, $false, $true, 1, $false, "32+19=51", 2) | Out-Null
$d.Content.Find.Execute("62-6=56", $true, $false, $false, $false, $false, $true, 1, $false, "61-29=32", 2) | Out-Null
$d.Content.Find.Execute("74-48=26", $true, $false, $false, $false, $false, $true, 1, $false, "44+18=62", 2) | Out-Null
$d.Content.Find.Execute("82-73=9", $true, $false, $false, $false, $false, $true, 1, $false, "25+68=93", 2) | Out-Null
$d.Content.Find.Execute("18+13=31", $true, $false, $false, $false, $false, $true, 1, $false, "32+19=51", 2) | Out-Null
$d.Content.Find.Execute("65+27=92", $true, $false, $false, $false, $false, $true, 1, $false, "54+29=83", 2) | Out-Null
$d.Content.Find.Execute("84+8=92", $true, $false, $false, $false, $false, $true, 1, $false, "83-35=48", 2) | Out-Null
$d.Content.Find.Execute("29+56=85", $true, $false, $false, $false, $false, $true, 1, $false, "91-19=72", 2) | Out-Null
$d.Content.Find.Execute("90-53=37", $true, $false, $false, $false, $false, $true, 1, $false, "58+26=84", 2) | Out-Null
$d.Content.Find.Execute("64+9=73", $true, $false, $false, $false, $false, $true, 1, $false, "39+12=51", 2) | Out-Null
$d.Content.Find.Execute("51-44=7", $true, $false, $false, $false, $false, $true, 1, $false, "15+47=62", 2) | Out-Null
$d.Content.Find.Execute("60-54=6", $true, $false, $false, $false, $false, $true, 1, $false, "71-44=27", 2) | Out-Null
$d.Content.Find.Execute("92-36=56", $true, $false, $false, $false, $false, $true, 1, $false, "50-47=3", 2) | Out-Null
$d.Content.Find.Execute("18+36=54", $true, $false, $false, $false, $false, $true, 1, $false, "72-33=39", 2) | Out-Null
$d.Content.Find.Execute("82-58=24", $true, $false, $false, $false, $false, $true, 1, $false, "6+29=35", 2) | Out-Null
$d.Content.Find.Execute("19+58=77", $true, $false, $false, $false, $false, $true, 1, $false, "14+58=72", 2) | Out-Null
$d.Content.Find.Execute("79+17=96", $true, $false, $false, $false, $false, $true, 1, $false, "9+25=34", 2) | Out-Null
$d.Content.Find.Execute("9+55=64", $true, $false, $false, $false, $false, $true, 1, $false, "50-12=38", 2) | Out-Null
$d.Content.Find.Execute("80-21=59", $true, $false, $false, $false, $false, $true, 1, $false, "17+34=51", 2) | Out-Null
$d.Content.Find.Execute("72-54=18", $true, $false, $false, $false, $false, $true, 1, $false, "39+52=91", 2) | Out-Null
$d.Content.Find.Execute("91-84=7", $true, $false, $false, $false, $false, $true, 1, $false, "21-2=19", 2) | Out-Null
$d.Content.Find.Execute("28+9=37", $true, $false, $false, $false, $false, $true, 1, $false, "18+18=36", 2) | Out-Null
$d.Content.Find.Execute("61-55=6", $true, $false, $false, $false, $false, $true, 1, $false, "38+54=92", 2) | Out-Null
$d.Content.Find.Execute("8+57=65", $true, $false, $false, $false, $false, $true, 1, $false, "94-57=37", 2) | Out-Null
$d.Content.Find.Execute("97-78=19", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=18", 2) | Out-Null
$d.Content.Find.Execute("57+28=85", $true, $false, $false, $false, $false, $true, 1, $false, "14-5=9", 2) | Out-Null
$d.Content.Find.Execute("27+18=45", $true, $false, $false, $false, $false, $true, 1, $false, "48-39=9", 2) | Out-Null
$d.Content.Find.Execute("19+22=41", $true, $false, $false, $false, $false, $true, 1, $false, "53-35=18", 2) | Out-Null
$d.Content.Find.Execute("81-77=4", $true, $false, $false, $false, $false, $true, 1, $false, "19+29=48", 2) | Out-Null
$d.Content.Find.Execute("17+29=46", $true, $false, $false, $false, $false, $true, 1, $false, "37+29=66", 2) | Out-Null
$d.Content.Find.Execute("23+28=51", $true, $false, $false, $false, $false, $true, 1, $false, "34+38=72", 2) | Out-Null
$d.Content.Find.Execute("47+47=94", $true, $false, $false, $false, $false, $true, 1, $false, "81-27=54", 2) | Out-Null
$d.Content.Find.Execute("41-29=12", $true, $false, $false, $false, $false, $true, 1, $false, "15+27=42", 2) | Out-Null
$d.Content.Find.Execute("4+49=53", $true, $false, $false, $false, $false, $true, 1, $false, "49+14=63", 2) | Out-Null
$d.Content.Find.Execute("50-32=18", $true, $false, $false, $false, $false, $true, 1, $false, "24+59=83", 2) | Out-Null
$d.Content.Find.Execute("3+89=92", $true, $false, $false, $false, $false, $true, 1, $false, "16+56=72", 2) | Out-Null
$d.Content.Find.Execute("29+25=54", $true, $false, $false, $false, $false, $true, 1, $false, "80-71=9", 2) | Out-Null
$d.Content.Find.Execute("80-4=76", $true, $false, $false, $false, $false, $true, 1, $false, "5+6=11", 2) | Out-Null
$d.Content.Find.Execute("58+9=67", $true, $false, $false, $false, $false, $true, 1, $false, "72-57=15", 2) | Out-Null
$d.Content.Find.Execute("15+79=94", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=62", 2) | Out-Null
$d.Content.Find.Execute("54+37=91", $true, $false, $false, $false, $false, $true, 1, $false, "49+45=94", 2) | Out-Null
$d.Content.Find.Execute("60-23=37", $true, $false, $false, $false, $false, $true, 1, $false, "54+27=81", 2) | Out-Null
$d.Content.Find.Execute("87-38=49", $true, $false, $false, $false, $false, $true, 1, $false, "13+18=31", 2) | Out-Null
$d.Content.Find.Execute("12+9=21", $true, $false, $false, $false, $false, $true, 1, $false, "30-29=1", 2) | Out-Null
$d.Content.Find.Execute("40-9=31", $true, $false, $false, $false, $false, $true, 1, $false, "94-7=87", 2) | Out-Null
$d.Content.Find.Execute("16+28=44", $true, $false, $false, $false, $false, $true, 1, $false, "53+28=81", 2) | Out-Null
$d.Content.Find.Execute("89+7=96", $true, $false, $false, $false, $false, $true, 1, $false, "66+19=85", 2) | Out-Null
$d.Content.Find.Execute("65+18=83", $true, $false, $false, $false, $false, $true, 1, $false, "70-37=33", 2) | Out-Null
$d.Content.Find.Execute("55-26=29", $true, $false, $false, $false, $false, $true, 1, $false, "64-9=55", 2) | Out-Null
$d.Content.Find.Execute("82-54=28", $true, $false, $false, $false, $false, $true, 1, $false, "31-13=18", 2) | Out-Null
$d.Content.Find.Execute("85-7=78", $true, $false, $false, $false, $false, $true, 1, $false, "21-17=4", 2) | Out-Null
$d.Content.Find.Execute("90-16=74", $true, $false, $false, $false, $false, $true, 1, $false, "9+35=44", 2) | Out-Null
$d.Content.Find.Execute("27+29=56", $true, $false, $false, $false, $false, $true, 1, $false, "13+59=72", 2) | Out-Null
$d.Content.Find.Execute("35-18=17", $true, $false, $false, $false, $false, $true, 1, $false, "24+67=91", 2) | Out-Null
$d.Content.Find.Execute("5+48=53", $true, $false, $false, $false, $false, $true, 1, $false, "6+59=65", 2) | Out-Null
$d.Content.Find.Execute("43+28=71", $true, $false, $false, $false, $false, $true, 1, $false, "14+68=82", 2) | Out-Null

Write-Output "Replacements applied: 101"
